$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Jira id) first so it becomes shared-string index 68,
# matching insertion order in the target workbook.
$ws.Range("B2").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "OPQA-1934||`nOPQA-1935&OPQA-3687||`nOPQA-4230||OPQA-4229||`nOPQA-4231||OPQA-4232||`nOPQA-4636||OPQA-4261||`nOPQA-4244||OPQA-4264||`nOPQA-4265||OPQA-4237||`nOPQA-4239||OPQA-4240||`nOPQA-4246||OPQA-4248||`nOPQA-4252"

# Column C (Description) second -> shared-string index 69.
# Base the cell format on an existing wrapped/bordered cell (C21),
# then darken the font color to FF333333 as in the target style.
$ws.Range("C21").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "Verify that Forgot your password? Link is clickable on NEON Landing page and End note landing page||Verify that the system is navigating to Forgot Password page or not|| after clicking on Forgot your password? Link&Verify that||the system should support a ENW password reset workflow with the following configurations||Verify that system should not inform user that entered email is not found.||Verify that user should be able to enter email address in Forgot password page.||Verify that  forget password service should send a forgot password email when the email entered is registered in the system||Verify that the platform password reset service should send a platform forget password email with branding that corresponds with the originating application as per wireframe||Verify that When the password reset token in the email is valid|| upon clicking the password reset link in the the platform forget password email|| the user shall be taken to the External Password Reset Page||Verify that External Password Reset Page should have a new password field where the user enters their new password.||Verify that when reset Password Token already used user should be taken to sign in screen||Verify that upon successful submission of a password change|| The user should receive a password change confirmation email to the user's primary email address with branding that corresponds with the application that the user completed the password change||Verify that the password change confirmation email should reference the fact that credentials are shared across all products.||Verify that when the password reset token in the email is expired or already used|| upon clicking the password reset link in the the platform forget password email|| the user should be taken to the External Invalid Password Reset Token Page.||Verify that the email address on the External Invalid Password Reset Token Page should be pre-populated with the email address that matches the email that the forgot password email was sent.||Verify that user who clicks the submit button on the the External Invalid Password Reset Token page|| should be taken to the target application sign in page.||Verify that when Email address is known from password reset token||error message 'The email address is prepopulated.' should be displayed and email address field should be editable||Verify that when Email address is not known from password reset token||email address field should be blank and user should be able to enter any email address||Verify that error message Please enter a valid email address.should be displayed in red color when user enters email address in wrong format"
$ws.Range("C23").Font.Color = 3355443

# Column A (TCID) last -> shared-string index 70.
$ws.Range("A2").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "IPAIAM100"

# Column D (Runmode) reuses the existing "Y" shared string.
$ws.Range("D2").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = "Y"

# Column E (Results) stays blank, matching formatting of E2.
$ws.Range("E2").Copy()
$ws.Range("E23").PasteSpecial(-4122)

# Match the tall row height used to show the full wrapped description.
$ws.Rows.Item(23).RowHeight = 409.5

# Select the new row, mirroring the saved view state in the target file.
$ws.Range("A23").Select() | Out-Null
